$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the revised figures on the existing last row (row 74) ---
$ws.Range("B74").Value = 33754
$ws.Range("E74").Value = 33474
$ws.Range("G74").Value = 1050

# --- Append the new quarterly row (row 75) ---
# Column A holds a "dd-mm-yyyy"-looking label that must stay plain text
# (shared string), not get auto-converted to a date serial. Writing it
# through a formula first (so Excel treats it as text, not a date) and
# then collapsing the formula down to its cached value with a
# copy / paste-values round trip keeps it a clean string cell with no
# extra number-format/style baggage.
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy() | Out-Null
$ws.Range("A75").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("B75").Value = 37860
$ws.Range("C75").Value = 142
$ws.Range("D75").Value = 142
$ws.Range("E75").Value = 37718
$ws.Range("F75").Value = 36692
$ws.Range("G75").Value = 1026
